# lab06 addition + fatal-error fixes
#
# Converts a target EMU value into the "points" value that must be
# assigned to a Shape's Left/Top/Width/Height (which PowerPoint's COM
# object model stores as a 32-bit Single, in points = EMU/12700) so
# that, after PowerPoint converts it back to EMU for storage in the
# OOXML, it lands exactly on the desired EMU value.
function EmuToPt($emu) {
    if ($emu -eq 0) { return 0.0 }
    $base = $emu / 12700.0
    for ($i = 0; $i -lt 20000; $i++) {
        $cand = $base + $i * 0.0000001
        $back = [int]([float]$cand * 12700.0)
        if ($back -eq $emu) {
            return $cand
        }
    }
    return $base
}

function Find-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        if ($shapes.Item($i).Id -eq $id) {
            return $shapes.Item($i)
        }
    }
    return $null
}

function Set-ShapeRectEmu($shape, $x, $y, $cx, $cy) {
    $shape.Left = EmuToPt $x
    $shape.Top = EmuToPt $y
    $shape.Width = EmuToPt $cx
    $shape.Height = EmuToPt $cy
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Fatal-error fix: re-point / resize the straight connector (old id
#    571, "直接连接符 570") that sits to the right of the instruction
#    register area.
# ---------------------------------------------------------------------
$conn571 = Find-ShapeById $s.Shapes 571
Set-ShapeRectEmu $conn571 8170286 546586 3434 2942212

# ---------------------------------------------------------------------
# 2) Fatal-error fix: the "Inst." label was split across two runs
#    ("I" + "nst."); merge it back into a single run. Force a real
#    text change (no-op Text assignments are not applied) before
#    setting the final value so the two runs actually collapse to one.
# ---------------------------------------------------------------------
$lbl257 = Find-ShapeById $s.Shapes 257
$lbl257.TextFrame.TextRange.Text = "__tmp__"
$lbl257.TextFrame.TextRange.Text = "Inst."

# ---------------------------------------------------------------------
# 3) Add part of lab06: a small AND/mux style doodle copied from the
#    existing matching shapes elsewhere on the slide (so the new
#    shapes pick up the exact same theme-based p:style block) then
#    repositioned, resized, relabeled to match the new diagram.
# ---------------------------------------------------------------------

# 3a) "椭圆 258" (Oval) -- clone of the last shape on the slide ("椭圆
#     257"), which already uses the identical dot style + size.
$ovalSrc = Find-ShapeById $s.Shapes 258
$oval259 = $ovalSrc.Duplicate()
Set-ShapeRectEmu $oval259 8147427 3114825 45719 45719
$oval259.Name = "椭圆 258"

# 3b) "直接连接符 259" (straight connector, no flip) -- clone of an
#     existing un-flipped connector that has the same dk1 style.
$connSrc260 = Find-ShapeById $s.Shapes 262
$conn260 = $connSrc260.Duplicate()
Set-ShapeRectEmu $conn260 5484354 3491690 2685932 0
$conn260.Name = "直接连接符 259"

# 3c) "圆角矩形 266" (rounded rectangle labelled "~") -- clone of the
#     existing "&" gate rectangle (same accent2 style + font), then
#     relabel its text.
$rrSrc267 = Find-ShapeById $s.Shapes 476
$rr267 = $rrSrc267.Duplicate()
Set-ShapeRectEmu $rr267 6559201 3403481 238496 164876
$rr267.Name = "圆角矩形 266"
$rr267.TextFrame.TextRange.Text = "~"

# 3d) "直接连接符 267" (straight connector, flipped vertically) --
#     clone of an existing flipV connector with the same dk1 style.
$connSrc268 = Find-ShapeById $s.Shapes 269
$conn268 = $connSrc268.Duplicate()
Set-ShapeRectEmu $conn268 5484354 3489730 1 198000
$conn268.Name = "直接连接符 267"

# 3e) "圆角矩形 265" (rounded rectangle labelled "&") -- clone of the
#     same "&" gate rectangle as 3c (keeps "&" text as-is).
$rrSrc266 = Find-ShapeById $s.Shapes 476
$rr266 = $rrSrc266.Duplicate()
Set-ShapeRectEmu $rr266 5300599 3586070 257883 164876
$rr266.Name = "圆角矩形 265"
$rr266.TextFrame.TextRange.Text = "X"
$rr266.TextFrame.TextRange.Text = "&"
